# Auto-generated Excel COM-interop script.
#
# The workbook tracks FFXIV crafting-leve profitability; columns H..N on
# each job sheet (currentAveragePrice / currentAveragePriceNQ / HQ,
# LevePriceNQ / HQ, LeveProfitNQ / HQ) are plain cached numbers refreshed
# from market-board data by a scheduled runner -- there are no formulas to
# recompute, so each touched cell is written with its new literal value.
# A handful of rows gain or lose a trailing M/N cell entirely (previously
# absent profit figures becoming computable, or vice versa); those use
# ClearContents()/a fresh .Value assignment to match.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 4141.0625
$ws.Range("I62").Value = 4236.4614
$ws.Range("J62").Value = 3727.6667
$ws.Range("K62").Value = 4236.4614
$ws.Range("L62").Value = 3727.6667
$ws.Range("M62").Value = -3612.4614
$ws.Range("N62").Value = -4975.6667
$ws.Range("H65").Value = 4141.0625
$ws.Range("I65").Value = 4236.4614
$ws.Range("J65").Value = 3727.6667
$ws.Range("K65").Value = 21182.307
$ws.Range("L65").Value = 18638.3335
$ws.Range("M65").Value = -18062.307
$ws.Range("N65").Value = -24878.3335
$ws.Range("H69").Value = 5984.5
$ws.Range("I69").Value = 5000
$ws.Range("K69").Value = 15000
$ws.Range("M69").Value = -14126
$ws.Range("H72").Value = 5984.5
$ws.Range("I72").Value = 5000
$ws.Range("K72").Value = 45000
$ws.Range("M72").Value = -40632
$ws.Range("H86").Value = 4175.077
$ws.Range("J86").Value = 3695.8333
$ws.Range("L86").Value = 3695.8333
$ws.Range("N86").Value = -5941.8333
$ws.Range("H89").Value = 4175.077
$ws.Range("J89").Value = 3695.8333
$ws.Range("L89").Value = 18479.1665
$ws.Range("N89").Value = -29711.1665
$ws.Range("H94").Value = 12999.6
$ws.Range("I94").Value = 8749.75
$ws.Range("K94").Value = 8749.75
$ws.Range("M94").Value = -8298.75
$ws.Range("H96").Value = 10000
$ws.Range("J96").Value = 10000
$ws.Range("L96").Value = 30000
$ws.Range("N96").Value = -32746
$ws.Range("H107").Value = 1487.1892
$ws.Range("I107").Value = 1099.16
$ws.Range("J107").Value = 2295.5833
$ws.Range("K107").Value = 1099.16
$ws.Range("L107").Value = 2295.5833
$ws.Range("M107").Value = 820.8399999999999
$ws.Range("N107").Value = -6135.5833
$ws.Range("H137").Value = 1931.75
$ws.Range("I137").Value = 1996.9231
$ws.Range("K137").Value = 5990.7693
$ws.Range("M137").Value = -3440.7693
$ws.Range("H141").Value = 9931.789000000001
$ws.Range("J141").Value = 66249.5
$ws.Range("L141").Value = 198748.5
$ws.Range("N141").Value = -209108.5

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3447.3044
$ws.Range("I32").Value = 2172.1035
$ws.Range("K32").Value = 2172.1035
$ws.Range("M32").Value = -1885.1035
$ws.Range("H45").Value = 1341.8235
$ws.Range("J45").Value = 1432.6428
$ws.Range("L45").Value = 1432.6428
$ws.Range("N45").Value = -2186.6428
$ws.Range("H63").Value = 1680.875
$ws.Range("I63").Value = 1564
$ws.Range("K63").Value = 1564
$ws.Range("M63").Value = -878
$ws.Range("H66").Value = 1680.875
$ws.Range("I66").Value = 1564
$ws.Range("K66").Value = 7820
$ws.Range("M66").Value = -4388
$ws.Range("H102").Value = 23834986
$ws.Range("I102").Value = 29413686
$ws.Range("J102").Value = 125512.25
$ws.Range("K102").Value = 29413686
$ws.Range("L102").Value = 125512.25
$ws.Range("M102").Value = -29412064
$ws.Range("N102").Value = -128756.25
$ws.Range("H122").Value = 10533565
$ws.Range("I122").Value = 11118407
$ws.Range("K122").Value = 33355221
$ws.Range("M122").Value = -33352771
$ws.Range("H132").Value = 10873.127
$ws.Range("I132").Value = 7319.7812
$ws.Range("J132").Value = 43360.855
$ws.Range("K132").Value = 21959.3436
$ws.Range("L132").Value = 130082.565
$ws.Range("M132").Value = -19429.3436
$ws.Range("N132").Value = -135142.565

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3165.647
$ws.Range("I86").Value = 3726.4
$ws.Range("K86").Value = 3726.4
$ws.Range("M86").Value = -2603.4
$ws.Range("H89").Value = 3165.647
$ws.Range("I89").Value = 3726.4
$ws.Range("K89").Value = 18632
$ws.Range("M89").Value = -13016
$ws.Range("H94").Value = 25643216
$ws.Range("J94").Value = 125005976
$ws.Range("L94").Value = 125005976
$ws.Range("N94").Value = -125006878
$ws.Range("H99").Value = 52632636
$ws.Range("I99").Value = 66667730
$ws.Range("K99").Value = 66667730
$ws.Range("M99").Value = -66666232
$ws.Range("H105").Value = 1287.375
$ws.Range("J105").Value = 455
$ws.Range("L105").Value = 455
$ws.Range("N105").Value = -3949
$ws.Range("H134").Value = 27780362
$ws.Range("I134").Value = 27780362
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 83341086
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -83338551
$ws.Range("N134").ClearContents()

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 956.619
$ws.Range("I16").Value = 777.1667
$ws.Range("K16").Value = 777.1667
$ws.Range("M16").Value = -490.1667
$ws.Range("H31").Value = 1914
$ws.Range("I31").Value = 1427.8
$ws.Range("J31").Value = 3129.5
$ws.Range("K31").Value = 1427.8
$ws.Range("L31").Value = 3129.5
$ws.Range("M31").Value = -1132.8
$ws.Range("N31").Value = -3719.5
$ws.Range("H34").Value = 1914
$ws.Range("I34").Value = 1427.8
$ws.Range("J34").Value = 3129.5
$ws.Range("K34").Value = 1427.8
$ws.Range("L34").Value = 3129.5
$ws.Range("M34").Value = -1225.8
$ws.Range("N34").Value = -3533.5
$ws.Range("H113").Value = 956.619
$ws.Range("I113").Value = 777.1667
$ws.Range("K113").Value = 777.1667
$ws.Range("M113").Value = 1392.8333
$ws.Range("H122").Value = 1945.5555
$ws.Range("I122").Value = 1963.8125
$ws.Range("K122").Value = 5891.4375
$ws.Range("M122").Value = -3441.4375
$ws.Range("H132").Value = 2512.25
$ws.Range("I132").Value = 2512.25
$ws.Range("K132").Value = 7536.75
$ws.Range("M132").Value = -5006.75

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H106").Value = 12932.667
$ws.Range("J106").Value = 13499
$ws.Range("L106").Value = 40497
$ws.Range("N106").Value = -42389
$ws.Range("H119").Value = 2209.3333
$ws.Range("I119").Value = 2209.3333
$ws.Range("K119").Value = 6627.999899999999
$ws.Range("M119").Value = -1789.999899999999

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5934.4443
$ws.Range("I80").Value = 6068.8335
$ws.Range("J80").Value = 5665.6665
$ws.Range("K80").Value = 6068.8335
$ws.Range("L80").Value = 5665.6665
$ws.Range("M80").Value = -5070.8335
$ws.Range("N80").Value = -7661.6665
$ws.Range("H83").Value = 5934.4443
$ws.Range("I83").Value = 6068.8335
$ws.Range("J83").Value = 5665.6665
$ws.Range("K83").Value = 30344.1675
$ws.Range("L83").Value = 28328.3325
$ws.Range("M83").Value = -25352.1675
$ws.Range("N83").Value = -38312.3325
$ws.Range("H96").Value = 30118.5
$ws.Range("I96").Value = 30237
$ws.Range("J96").Value = 30000
$ws.Range("K96").Value = 30237
$ws.Range("L96").Value = 30000
$ws.Range("M96").Value = -27491
$ws.Range("N96").Value = -35492
$ws.Range("H102").Value = 927.78125
$ws.Range("I102").Value = 951.25806
$ws.Range("K102").Value = 951.25806
$ws.Range("M102").Value = 670.74194
$ws.Range("H122").Value = 34484276
$ws.Range("I122").Value = 41667900
$ws.Range("K122").Value = 125003700
$ws.Range("M122").Value = -125001250
$ws.Range("H126").Value = 2558.2222
$ws.Range("I126").Value = 2739.4285
$ws.Range("J126").Value = 1924
$ws.Range("K126").Value = 8218.2855
$ws.Range("L126").Value = 5772
$ws.Range("M126").Value = -5748.2855
$ws.Range("N126").Value = -10712
$ws.Range("H132").Value = 951.1111
$ws.Range("I132").Value = 951.1111
$ws.Range("K132").Value = 2853.3333
$ws.Range("M132").Value = -323.3332999999998

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6171.5835
$ws.Range("I7").Value = 4966.857
$ws.Range("J7").Value = 7858.2
$ws.Range("K7").Value = 4966.857
$ws.Range("L7").Value = 7858.2
$ws.Range("M7").Value = -4854.857
$ws.Range("N7").Value = -8082.2
$ws.Range("H19").Value = 10251.5
$ws.Range("I19").Value = 10251.5
$ws.Range("K19").Value = 10251.5
$ws.Range("M19").Value = -10081.5
$ws.Range("H22").Value = 1891.1538
$ws.Range("J22").Value = 1841.1818
$ws.Range("L22").Value = 1841.1818
$ws.Range("N22").Value = -2431.1818
$ws.Range("H27").Value = 1891.1538
$ws.Range("J27").Value = 1841.1818
$ws.Range("L27").Value = 1841.1818
$ws.Range("N27").Value = -2055.1818
$ws.Range("H40").Value = 5212.25
$ws.Range("I40").Value = 5012.4346
$ws.Range("K40").Value = 5012.4346
$ws.Range("M40").Value = -4876.4346
$ws.Range("H55").Value = 1956.92
$ws.Range("I55").Value = 453.35715
$ws.Range("J55").Value = 3870.5454
$ws.Range("K55").Value = 453.35715
$ws.Range("L55").Value = 3870.5454
$ws.Range("M55").Value = -280.35715
$ws.Range("N55").Value = -4216.5454
$ws.Range("H100").Value = 4279.4
$ws.Range("J100").Value = 4865.5
$ws.Range("L100").Value = 4865.5
$ws.Range("N100").Value = -5947.5
$ws.Range("H126").Value = 6171.5835
$ws.Range("I126").Value = 4966.857
$ws.Range("J126").Value = 7858.2
$ws.Range("K126").Value = 14900.571
$ws.Range("L126").Value = 23574.6
$ws.Range("M126").Value = -12430.571
$ws.Range("N126").Value = -28514.6

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3450.9167
$ws.Range("I81").Value = 1355.1111
$ws.Range("J81").Value = 9738.333000000001
$ws.Range("K81").Value = 2710.2222
$ws.Range("L81").Value = 19476.666
$ws.Range("M81").Value = -1649.2222
$ws.Range("N81").Value = -21598.666
$ws.Range("H84").Value = 3450.9167
$ws.Range("I84").Value = 1355.1111
$ws.Range("J84").Value = 9738.333000000001
$ws.Range("K84").Value = 13551.111
$ws.Range("L84").Value = 97383.33
$ws.Range("M84").Value = -8247.111000000001
$ws.Range("N84").Value = -107991.33
